$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("L2").Value = 133.36
$ws1.Range("L30").Value = "2 de 28"

$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F2").Value = 133.36
$ws2.Range("F29").Value = 27.73
$ws2.Range("F30").Value = 983.25

$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D10").Value = 27.73
$ws3.Range("E10").Value = 622.52
$ws3.Range("F10").Value = 0.04264513648596693

$ws3.Range("D16").Value = 955.52
$ws3.Range("E16").Value = 17843.09
$ws3.Range("F16").Value = 0.05082929003793366

$ws3.Range("D19").Value = 983.25
$ws3.Range("E19").Value = 28554.54107555787
$ws3.Range("F19").Value = 0.03328786494172294
